# "Generate Report for Handoff"
#
# The localization-status workbook tracks, per source file, the most recent
# handoff to each locale (zh-cn / de-de) plus a rolled-up "Overview" sheet.
# A new handoff just completed for row 5 (b000d19f-71b8-4277-8b7e-804a03fbe4f2),
# so its "Latest Handoff Datetime" timestamps move forward on the locale
# sheets, and the "Overview" sheet's "Latest Handoff Date" for that row
# (the max across locales) advances to match the de-de value.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn: row 5 ("b000d19f-...md") gets a fresh handoff timestamp.
$zhcn.Range("E5").Value = "2016-03-21 03:57:14"

# de-de: row 5 ("b000d19f-...md") gets a fresh handoff timestamp.
$dede.Range("E5").Value = "2016-03-21 03:57:23"

# Overview: row 5's "Latest Handoff Date" reflects the newest handoff
# across locales (de-de's 03:57:23).
$overview.Range("D5").Value = "2016-03-21 03:57:23"
